# Updates the date header and the multiplication problems in the table.
$d = $word.ActiveDocument

# Map of old text -> new text (each old value is unique in the document).
$replacements = @(
    @("2026-01-17 Saturday", "2026-01-18 Sunday"),
    @("136×3=", "953×7="),
    @("285×7=", "409×3="),
    @("516×8=", "144×2="),
    @("435×5=", "849×4="),
    @("308×4=", "398×3="),
    @("870×9=", "181×2="),
    @("931×6=", "220×3="),
    @("951×9=", "203×4="),
    @("363×5=", "401×9="),
    @("780×9=", "155×4="),
    @("624×4=", "308×7="),
    @("361×8=", "482×5="),
    @("557×5=", "477×6="),
    @("260×3=", "322×3="),
    @("729×3=", "224×4="),
    @("855×9=", "897×7="),
    @("929×2=", "306×2="),
    @("202×4=", "810×7="),
    @("144×4=", "282×5="),
    @("562×7=", "508×2="),
    @("716×5=", "836×7="),
    @("850×6=", "352×9="),
    @("233×5=", "383×4="),
    @("746×9=", "236×5="),
    @("216×9=", "369×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
